$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 16920316
$summary.Range("B7").Value = 1666320
$summary.Range("B8").Value = 1080

# --- Purchasing Costs sheet ---
$pc = $wb.Worksheets.Item("Purchasing Costs")
for ($r = 2; $r -le 11; $r++) {
    $pc.Cells.Item($r, 8).Value = 166632   # column H: Total Holding Cost
    $pc.Cells.Item($r, 9).Value = 108      # column I: Total Fixed Cost
}
